$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a numeric-looking value that must be stored as text,
# so pre-format the cell as Text before entering the value.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "123456789"
$ws.Range("B1").Value = "raful9"
$ws.Range("C1").Value = "123456a!"
$ws.Range("D1").Value = "rafael"
